$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 18041.666
$ws.Range("J123").Value = 18041.666
$ws.Range("L123").Value = 18041.666
$ws.Range("N123").Value = -27841.666
$ws.Range("H132").Value = 901.3158
$ws.Range("I132").Value = 779.96075
$ws.Range("J132").Value = 1932.8334
$ws.Range("K132").Value = 2339.88225
$ws.Range("L132").Value = 5798.5002
$ws.Range("M132").Value = 190.1177500000003
$ws.Range("N132").Value = -10858.5002
$ws.Range("H135").Value = 2052.262
$ws.Range("I135").Value = 1547.6052
$ws.Range("J135").Value = 6846.5
$ws.Range("K135").Value = 13928.4468
$ws.Range("L135").Value = 61618.5
$ws.Range("M135").Value = -11393.4468
$ws.Range("N135").Value = -66688.5
$ws.Range("H137").Value = 1694.075
$ws.Range("I137").Value = 1435.2333
$ws.Range("K137").Value = 4305.699900000001
$ws.Range("M137").Value = -1755.699900000001
$ws.Range("H138").Value = 1223.98
$ws.Range("I138").Value = 694.5854
$ws.Range("J138").Value = 1591.8644
$ws.Range("K138").Value = 2083.7562
$ws.Range("L138").Value = 4775.593199999999
$ws.Range("M138").Value = 3056.2438
$ws.Range("N138").Value = -15055.5932

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6785.2856
$ws.Range("I32").Value = 5543.436
$ws.Range("J32").Value = 11628.5
$ws.Range("K32").Value = 5543.436
$ws.Range("L32").Value = 11628.5
$ws.Range("M32").Value = -5256.436
$ws.Range("N32").Value = -12202.5
$ws.Range("H61").Value = 4748.121
$ws.Range("I61").Value = 5426.4
$ws.Range("J61").Value = 2628.5
$ws.Range("K61").Value = 5426.4
$ws.Range("L61").Value = 2628.5
$ws.Range("M61").Value = -5214.4
$ws.Range("N61").Value = -3052.5
$ws.Range("H136").Value = 4748.121
$ws.Range("I136").Value = 5426.4
$ws.Range("J136").Value = 2628.5
$ws.Range("K136").Value = 16279.2
$ws.Range("L136").Value = 7885.5
$ws.Range("M136").Value = -13729.2
$ws.Range("N136").Value = -12985.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1831.1875
$ws.Range("I86").Value = 1780
$ws.Range("J86").Value = 1916.5
$ws.Range("K86").Value = 1780
$ws.Range("L86").Value = 1916.5
$ws.Range("M86").Value = -657
$ws.Range("N86").Value = -4162.5
$ws.Range("H89").Value = 1831.1875
$ws.Range("I89").Value = 1780
$ws.Range("J89").Value = 1916.5
$ws.Range("K89").Value = 8900
$ws.Range("L89").Value = 9582.5
$ws.Range("M89").Value = -3284
$ws.Range("N89").Value = -20814.5
$ws.Range("H94").Value = 1973.9375
$ws.Range("I94").Value = 863.8333
$ws.Range("J94").Value = 2640
$ws.Range("K94").Value = 863.8333
$ws.Range("L94").Value = 2640
$ws.Range("M94").Value = -412.8333
$ws.Range("N94").Value = -3542
$ws.Range("H99").Value = 58824852
$ws.Range("I99").Value = 76924110
$ws.Range("J99").Value = 2249.75
$ws.Range("K99").Value = 76924110
$ws.Range("L99").Value = 2249.75
$ws.Range("M99").Value = -76922612
$ws.Range("N99").Value = -5245.75
$ws.Range("H107").Value = 804.96155
$ws.Range("I107").Value = 634.1177
$ws.Range("J107").Value = 1127.6666
$ws.Range("K107").Value = 634.1177
$ws.Range("L107").Value = 1127.6666
$ws.Range("M107").Value = 1285.8823
$ws.Range("N107").Value = -4967.6666
$ws.Range("H134").Value = 4335.255
$ws.Range("I134").Value = 6444.091
$ws.Range("J134").Value = 2735.4482
$ws.Range("K134").Value = 19332.273
$ws.Range("L134").Value = 8206.3446
$ws.Range("M134").Value = -16797.273
$ws.Range("N134").Value = -13276.3446

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1269.6842
$ws.Range("I58").Value = 679
$ws.Range("J58").Value = 2282.2856
$ws.Range("K58").Value = 679
$ws.Range("L58").Value = 2282.2856
$ws.Range("M58").Value = -476
$ws.Range("N58").Value = -2688.2856
$ws.Range("H107").Value = 10753466
$ws.Range("I107").Value = 13889463
$ws.Range("J107").Value = 1477.5714
$ws.Range("K107").Value = 13889463
$ws.Range("L107").Value = 1477.5714
$ws.Range("M107").Value = -13887543
$ws.Range("N107").Value = -5317.5714
$ws.Range("H132").Value = 2066.8572
$ws.Range("I132").Value = 1549.4166
$ws.Range("K132").Value = 4648.2498
$ws.Range("M132").Value = -2118.2498
$ws.Range("H134").Value = 1801.5397
$ws.Range("I134").Value = 2096.675
$ws.Range("J134").Value = 1288.2609
$ws.Range("K134").Value = 6290.025000000001
$ws.Range("L134").Value = 3864.7827
$ws.Range("M134").Value = -3755.025000000001
$ws.Range("N134").Value = -8934.7827
$ws.Range("H136").Value = 1269.6842
$ws.Range("I136").Value = 679
$ws.Range("J136").Value = 2282.2856
$ws.Range("K136").Value = 2037
$ws.Range("L136").Value = 6846.8568
$ws.Range("M136").Value = 513
$ws.Range("N136").Value = -11946.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 333334850
$ws.Range("J9").Value = 333334850
$ws.Range("L9").Value = 1000004550
$ws.Range("N9").Value = -1000004998
$ws.Range("H107").Value = 393.33334
$ws.Range("I107").Value = 239
$ws.Range("J107").Value = 434.8846
$ws.Range("K107").Value = 717
$ws.Range("L107").Value = 1304.6538
$ws.Range("M107").Value = 1203
$ws.Range("N107").Value = -5144.6538
$ws.Range("H131").Value = 1961614.5
$ws.Range("I131").Value = 7692591.5
$ws.Range("J131").Value = 1017.1842
$ws.Range("K131").Value = 23077774.5
$ws.Range("L131").Value = 3051.5526
$ws.Range("M131").Value = -23072734.5
$ws.Range("N131").Value = -13131.5526
$ws.Range("H132").Value = 8172914.5
$ws.Range("I132").Value = 3424.75
$ws.Range("J132").Value = 10896078
$ws.Range("K132").Value = 30822.75
$ws.Range("L132").Value = 98064702
$ws.Range("M132").Value = -28292.75
$ws.Range("N132").Value = -98069762
$ws.Range("H138").Value = 7248
$ws.Range("I138").Value = 7691.8423
$ws.Range("J138").Value = 3031.5
$ws.Range("K138").Value = 23075.5269
$ws.Range("L138").Value = 9094.5
$ws.Range("M138").Value = -17935.5269
$ws.Range("N138").Value = -19374.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 30145
$ws.Range("J32").Value = 30145
$ws.Range("L32").Value = 30145
$ws.Range("N32").Value = -30737
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 30000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 30000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -29515
$ws.Range("N42").ClearContents()
$ws.Range("H70").Value = 5603.0884
$ws.Range("I70").Value = 5748.3076
$ws.Range("J70").Value = 5131.125
$ws.Range("K70").Value = 5748.3076
$ws.Range("L70").Value = 5131.125
$ws.Range("M70").Value = -5478.3076
$ws.Range("N70").Value = -5671.125
$ws.Range("H73").Value = 5603.0884
$ws.Range("I73").Value = 5748.3076
$ws.Range("J73").Value = 5131.125
$ws.Range("K73").Value = 5748.3076
$ws.Range("L73").Value = 5131.125
$ws.Range("M73").Value = -4812.3076
$ws.Range("N73").Value = -7003.125
$ws.Range("H115").Value = 30000
$ws.Range("I115").Value = 30000
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 30000
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -28825
$ws.Range("N115").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 19608820
$ws.Range("I46").Value = 47619692
$ws.Range("J46").Value = 1210
$ws.Range("K46").Value = 47619692
$ws.Range("L46").Value = 1210
$ws.Range("M46").Value = -47619504
$ws.Range("N46").Value = -1586
$ws.Range("H132").Value = 17816704
$ws.Range("I132").Value = 21378928
$ws.Range("J132").Value = 5579.6
$ws.Range("K132").Value = 64136784
$ws.Range("L132").Value = 16738.8
$ws.Range("M132").Value = -64134254
$ws.Range("N132").Value = -21798.8
$ws.Range("H136").Value = 10021.823
$ws.Range("I136").Value = 7084.478
$ws.Range("K136").Value = 21253.434
$ws.Range("M136").Value = -18703.434

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H113").Value = 760.2182
$ws.Range("I113").Value = 668.57574
$ws.Range("J113").Value = 897.6818
$ws.Range("K113").Value = 2005.72722
$ws.Range("L113").Value = 2693.0454
$ws.Range("M113").Value = 164.27278
$ws.Range("N113").Value = -7033.0454
$ws.Range("H115").Value = 21000
$ws.Range("J115").Value = 21000
$ws.Range("L115").Value = 21000
$ws.Range("N115").Value = -24134
$ws.Range("H132").Value = 13417.456
$ws.Range("I132").Value = 15501.299
$ws.Range("J132").Value = 1782.6666
$ws.Range("K132").Value = 46503.897
$ws.Range("L132").Value = 5347.9998
$ws.Range("M132").Value = -43973.897
$ws.Range("N132").Value = -10407.9998
$ws.Range("H136").Value = 9618371
$ws.Range("I136").Value = 3633.0344
$ws.Range("J136").Value = 21741300
$ws.Range("K136").Value = 10899.1032
$ws.Range("L136").Value = 65223900
$ws.Range("M136").Value = -8349.1032
$ws.Range("N136").Value = -65229000
